# Apply cryptos list update (prices + volume%) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083  # subscript 3, used in PEPE price (0.0<sub>3</sub>0999)
$sub6 = [char]0x2086  # subscript 6, used in BabyDogeCoin price (0.0<sub>6</sub>0298)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.252.27'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.621.09'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.76'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.18'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.57%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.620.15'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.140'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.99%  '
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.368'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.41'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.100.84'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000181'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.356.73'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.625.51'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.77'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.81'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '356.55'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.30'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.67'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.61'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.03%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -4.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '69.61'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.757.95'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0${sub3}0999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '547.49'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.10'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("E33").Value = '  -3.92%  '
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.134'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -4.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.77'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.94'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.74%  '
$ws.Range("E40").Value = '  -2.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.26'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.16'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.10%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("E45").Value = '  -4.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0${sub6}0298"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.580'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '151.33'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.78'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.72'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.75%  '
$ws.Range("E51").Value = '  -1.74%  '
